# Update the player roster table on Sheet1.
# Columns: A = Oyuncu Adı (Player Name), B = Pozisyon (Position), C = Takım (Team)
# The table is re-sorted/updated and one extra row is added (Grayson Allen),
# while the previous "Bennedict Mathurin" row is replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Chris Paul",        "PG",        "San Antonio Spurs"),
    @("Dejounte Murray",   "PG,SG",     "New Orleans Pelicans"),
    @("Russell Westbrook", "PG,SG",     "Denver Nuggets"),
    @("Keon Johnson",      "PG,SG",     "Brooklyn Nets"),
    @("Khris Middleton",   "SF",        "Milwaukee Bucks"),
    @("Paolo Banchero",    "SF,PF",     "Orlando Magic"),
    @("Jaylen Brown",      "SG,SF",     "Boston Celtics"),
    @("Rudy Gobert",       "C",         "Minnesota Timberwolves"),
    @("Jakob Poeltl",      "C",         "Toronto Raptors"),
    @("Nikola Jokic",      "C",         "Denver Nuggets"),
    @("Jalen Green",       "PG,SG",     "Houston Rockets"),
    @("Buddy Hield",       "SG,SF",     "Golden State Warriors"),
    @("Caris LeVert",      "SG,SF",     "Cleveland Cavaliers"),
    @("Grayson Allen",     "PG,SG,SF",  "Phoenix Suns"),
    @("Pascal Siakam",     "SF,PF,C",   "Indiana Pacers"),
    @("Chet Holmgren",     "PF,C",      "Oklahoma City Thunder"),
    @("Jalen Suggs",       "PG,SG",     "Orlando Magic"),
    @("Deni Avdija",       "SF,PF",     "Portland Trail Blazers")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}
